# ------------------------------------------------------------------
# Re-theme "The Fabric of the Universe" (physics) into
# "The Enchanting World of Chemistry" (chemistry), per commit diff.
# ------------------------------------------------------------------
$d = $word.ActiveDocument

# Helper: locate a unique literal substring anywhere in the document
# and overwrite it with new text. We locate via Find (no replace
# argument, so Word's "smart quote" substitution on the replacement
# text never kicks in) and then assign straight onto the matched
# Range, which keeps straight apostrophes/quotes intact.
function Set-FoundText([string]$find, [string]$newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $find"
    }
    $rng.Text = $newText
}

# --- Title -----------------------------------------------------------
Set-FoundText "The Fabric of the Universe-  A Physicist's Perspective" `
              "The Enchanting World of Chemistry: Unveiling Nature's Symphony of Elements"

# --- Author name: "Isaac Newton" -> "Dr. Emily Carter" ---------------
Set-FoundText "Isaac Newton" "Dr. Emily Carter"

# --- Author e-mail paragraph ------------------------------------------
# Original runs: "isaac" / "." / "newton@physics" / "." / "com"
# Target text:    "ecarter@highschoolchem.edu"
$emailPara = $d.Paragraphs.Item(3)
$emailRange = $d.Range($emailPara.Range.Start, $emailPara.Range.End - 1)
$emailRange.Text = "ecarter@highschoolchem.edu"

# --- Body paragraph 1 (four "sentence" runs, separated by <br/><br/>) -
Set-FoundText "The cosmos, a symphony of celestial bodies, captivates our imagination and fuels scientific inquiry" `
              "In the vast tapestry of sciences, chemistry stands out as a symphony of elements, a harmonious dance of molecules, and a vibrant narrative of matter"

Set-FoundText " Physics unlocks the secrets of the universe, delving into the fundamental laws that govern matter and energy, space and time" `
              " It delves into the intricate relationships between substances, unveiling the secrets of their composition, properties, and interactions"

Set-FoundText " One of the most prominent theories in physics, quantum mechanics, has profoundly shaped our understanding of the universe at its smallest scales" `
              " Chemistry is the language of the natural world, spoken in the patterns of atoms, the eloquence of chemical bonds, and the ever-changing states of matter"

Set-FoundText " This fascinating realm, ruled by enigmatic particles and forces, has unveiled a hidden tapestry of quantum superposition, entanglement, and wave-particle duality, forever altering our perception of reality" `
              " Through chemistry, we gain the power to decipher nature's enigmatic whispers, unravel the complexities of our world, and harness its boundless potential for innovation"

Set-FoundText "As we peer deeper into the vastness of the cosmos, a symphony of celestial bodies reveals the profound unity underlying the universe's diverse phenomena" `
              "With each new discovery, chemistry opens doors to uncharted territories of knowledge, revealing the hidden wonders of the universe"

Set-FoundText " Einstein's theory of general relativity, a captivating tapestry of spacetime curvature, gravity, and the cosmic dance of celestial objects, has transformed our understanding of gravity and the cosmos" `
              " It empowers us to decode the intricate workings of life, unlock the secrets of disease, and devise ingenious solutions to global challenges"

Set-FoundText " General relativity's elegance and predictive power have enabled us to unlock the mysteries of black holes, gravitational waves, and the expansion of the universe, offering a glimpse into the breathtaking vastness of existence" `
              " Chemistry is the key to unraveling the mysteries of the cosmos, understanding the intricate mechanisms of our bodies, and developing revolutionary technologies that shape our future"

Set-FoundText "Delving into the realm of subatomic particles, the Standard Model of Physics gracefully orchestrates the intricate ballet of fundamental forces and particles, providing a comprehensive framework that encompasses the electromagnetic, weak, and strong interactions" `
              "As we delve deeper into the enchanting world of chemistry, we embark on a journey of exploration, experimentation, and enlightenment"

Set-FoundText " This symphony of subatomic interactions forms the foundation of matter and energy, dictating the properties of atoms, molecules, and the world we experience" `
              " We become alchemists, transforming ordinary substances into extraordinary materials, unlocking the secrets of chemical reactions, and witnessing the magic of transformations"

Set-FoundText " From the birth of stars to the fusion within them and the radiant melodies of atomic transitions, physics unlocks the secrets of energy transformation and the symphony of the universe" `
              " Chemistry invites us to explore the boundless possibilities of matter, to create new substances, and to understand the interconnectedness of all things"

# --- "Summary" heading is unchanged ------------------------------------

# --- Summary body paragraph --------------------------------------------
Set-FoundText "Physics, like a master conductor, weaves an intricate tapestry of knowledge, harmonizing universal laws, quantum mysteries, and subatomic dances" `
              "The realm of chemistry is an enchanting tapestry of elements, molecules, and interactions, offering a symphony of knowledge about the composition, properties, and transformations of matter"

Set-FoundText " From celestial symphonies in the cosmos to the enigmatic world of quantum particles, physics unveils the profound beauty and unity that underlies all of existence" `
              " It empowers us to understand the natural world, decode the complexities of life, and devise innovative solutions to global challenges"

# The final two sentences (and the "." run that used to separate them)
# collapse into a single continuous span of new text, still followed
# by the paragraph's closing "." run.
Set-FoundText ". This journey of exploration has not only enriched our understanding of the universe but has also propelled us to the forefront of technological advancements, enriching our lives in myriad ways. Thus, we embark on an endless quest to unravel the symphony of the universe, our perception of reality forever transformed by the brilliance of physics" `
              ". As we delve deeper into the enchanting world of chemistry, we embark on a journey of exploration, experimentation, and enlightenment, unlocking the secrets of matter and shaping our future"

# --- Trailing empty paragraph, newly added at the end of the body -----
$d.Paragraphs.Last.Range.InsertParagraphAfter()
